$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D10").Value = 4374
$ws.Range("H10").Value = 170

$ws.Range("H11").Select()
